$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 287, shifting the existing rows 287-324 down to 288-325.
$ws.Rows.Item(287).Insert()

# Populate the newly inserted row 287 with the new weekly record.
$ws.Range("A287").Value2 = 6
$ws.Range("B287").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C287").Value2 = "Metropolitana"
$ws.Range("D287").Value2 = 45127
$ws.Range("E287").Value2 = 13
$ws.Range("F287").Value2 = 100112022
$ws.Range("G287").Value2 = "Arveja Verde"
$ws.Range("H287").Value2 = "Perfection"
$ws.Range("I287").Value2 = "Primera"
$ws.Range("J287").Value2 = 350
$ws.Range("K287").Value2 = 25000
$ws.Range("L287").Value2 = 27000
$ws.Range("M287").Value2 = 26314
$ws.Range("N287").Value2 = "$/malla 25 kilos"
$ws.Range("O287").Value2 = "Provincia de Limarí"
$ws.Range("P287").Value2 = 1053
$ws.Range("Q287").Value2 = 25
$ws.Range("R287").Value2 = "Hortaliza"
